$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.116251945495605
$ws.Range("B1").Value = 1.8089519739151
$ws.Range("D1").Value = 1.827853918075562
$ws.Range("E1").Value = 1.109164834022522
